# Generate Report for Handoff
# The ffa9b05e-361c-4efb-966a-0babefbfbb56 file has moved from
# "Handed back: in sync with en-US" to "Ready for handoff", with fresh
# handoff timestamps recorded on the Overview sheet and the per-locale
# (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the ffa9b05e... file ---
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-33-11 14:33:07"

# --- zh-cn detail sheet: row 3 is the ffa9b05e... file ---
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-11 14:33:03"

# --- de-de detail sheet: row 3 is the ffa9b05e... file ---
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-11 14:33:07"
